# Auto-generated script applying market-data refresh changes to Maduin_Profits workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 2528.4614
$ws.Range("J88").Value = 2871.5454
$ws.Range("L88").Value = 2871.5454
$ws.Range("N88").Value = -3683.5454

$ws.Range("H91").Value = 2528.4614
$ws.Range("J91").Value = 2871.5454
$ws.Range("L91").Value = 2871.5454
$ws.Range("N91").Value = -5679.5454

$ws.Range("H111").Value = 1831.5
$ws.Range("I111").Value = 1913
$ws.Range("J111").Value = 1750
$ws.Range("K111").Value = 5739
$ws.Range("L111").Value = 5250
$ws.Range("M111").Value = -2672
$ws.Range("N111").Value = -11384

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5124
$ws.Range("I61").Value = 5500
$ws.Range("K61").Value = 5500
$ws.Range("M61").Value = -5288

$ws.Range("H74").Value = 747.8570999999999
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 747.8570999999999
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

$ws.Range("H132").Value = 1128.8462
$ws.Range("I132").Value = 1128.8462
$ws.Range("K132").Value = 3386.5386
$ws.Range("M132").Value = -856.5385999999999

$ws.Range("H136").Value = 5124
$ws.Range("I136").Value = 5500
$ws.Range("K136").Value = 16500
$ws.Range("M136").Value = -13950

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4182.25
$ws.Range("I105").Value = 3576.3333
$ws.Range("K105").Value = 3576.3333
$ws.Range("M105").Value = -1829.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 5644.8
$ws.Range("I58").Value = 2056
$ws.Range("K58").Value = 2056
$ws.Range("M58").Value = -1853

$ws.Range("H69").Value = 5000
$ws.Range("I69").Value = 5000
$ws.Range("K69").Value = 5000
$ws.Range("M69").Value = -4251

$ws.Range("H72").Value = 5000
$ws.Range("I72").Value = 5000
$ws.Range("K72").Value = 15000
$ws.Range("M72").Value = -11256

$ws.Range("H99").Value = 3608.6785
$ws.Range("I99").Value = 3062
$ws.Range("J99").Value = 4975.375
$ws.Range("K99").Value = 3062
$ws.Range("L99").Value = 4975.375
$ws.Range("M99").Value = -1564
$ws.Range("N99").Value = -7971.375

$ws.Range("H105").Value = 1506.6923
$ws.Range("I105").Value = 922
$ws.Range("J105").Value = 2442.2
$ws.Range("K105").Value = 922
$ws.Range("L105").Value = 2442.2
$ws.Range("M105").Value = 825
$ws.Range("N105").Value = -5936.2

$ws.Range("H126").Value = 3608.6785
$ws.Range("I126").Value = 3062
$ws.Range("J126").Value = 4975.375
$ws.Range("K126").Value = 9186
$ws.Range("L126").Value = 14926.125
$ws.Range("M126").Value = -6716
$ws.Range("N126").Value = -19866.125

$ws.Range("H132").Value = 3424
$ws.Range("I132").Value = 3746.5
$ws.Range("J132").Value = 1327.75
$ws.Range("K132").Value = 11239.5
$ws.Range("L132").Value = 3983.25
$ws.Range("M132").Value = -8709.5
$ws.Range("N132").Value = -9043.25

$ws.Range("H134").Value = 2199.8
$ws.Range("I134").Value = 2000
$ws.Range("K134").Value = 6000
$ws.Range("M134").Value = -3465

$ws.Range("H136").Value = 5644.8
$ws.Range("I136").Value = 2056
$ws.Range("K136").Value = 6168
$ws.Range("M136").Value = -3618

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 932.6
$ws.Range("I5").Value = 962.6667
$ws.Range("K5").Value = 2888.0001
$ws.Range("M5").Value = -2776.0001

$ws.Range("H12").Value = 756.3333
$ws.Range("J12").Value = 1224
$ws.Range("L12").Value = 3672
$ws.Range("N12").Value = -4018

$ws.Range("H38").Value = 1735
$ws.Range("I38").Value = 5012.5
$ws.Range("J38").Value = 96.25
$ws.Range("K38").Value = 15037.5
$ws.Range("L38").Value = 288.75
$ws.Range("M38").Value = -14690.5
$ws.Range("N38").Value = -982.75

$ws.Range("H118").Value = 550
$ws.Range("I118").Value = 550
$ws.Range("K118").Value = 1650
$ws.Range("M118").Value = -407

$ws.Range("H121").Value = 1741.5
$ws.Range("I121").Value = 833.3333
$ws.Range("J121").Value = 2286.4
$ws.Range("K121").Value = 2499.9999
$ws.Range("L121").Value = 6859.200000000001
$ws.Range("M121").Value = -1189.9999
$ws.Range("N121").Value = -9479.200000000001

$ws.Range("H134").Value = 872.5
$ws.Range("I134").Value = 872.5
$ws.Range("K134").Value = 2617.5
$ws.Range("M134").Value = 2452.5

$ws.Range("H135").Value = 932.6
$ws.Range("I135").Value = 962.6667
$ws.Range("K135").Value = 8664.0003
$ws.Range("M135").Value = -6129.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 49205.855
$ws.Range("J24").Value = 49205.855
$ws.Range("L24").Value = 49205.855
$ws.Range("N24").Value = -49551.855

$ws.Range("H113").Value = 1467.8334
$ws.Range("I113").Value = 1467.8334
$ws.Range("K113").Value = 1467.8334
$ws.Range("M113").Value = 702.1666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2879.75
$ws.Range("I61").Value = 2755.4285
$ws.Range("K61").Value = 2755.4285
$ws.Range("M61").Value = -2553.4285

$ws.Range("H106").Value = 29452.25
$ws.Range("J106").Value = 29452.25
$ws.Range("L106").Value = 29452.25
$ws.Range("N106").Value = -31976.25

$ws.Range("H113").Value = 2879.75
$ws.Range("I113").Value = 2755.4285
$ws.Range("K113").Value = 2755.4285
$ws.Range("M113").Value = -585.4285

$ws.Range("H136").Value = 4104.75
$ws.Range("I136").Value = 4035.7
$ws.Range("J136").Value = 4450
$ws.Range("K136").Value = 12107.1
$ws.Range("L136").Value = 13350
$ws.Range("M136").Value = -9557.099999999999
$ws.Range("N136").Value = -18450

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H76").Value = 20173
$ws.Range("J76").Value = 20173
$ws.Range("L76").Value = 20173
$ws.Range("N76").Value = -20803

$ws.Range("H79").Value = 20173
$ws.Range("J79").Value = 20173
$ws.Range("L79").Value = 20173
$ws.Range("N79").Value = -22357

$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws.Range("H132").Value = 1874.375
$ws.Range("I132").Value = 1856.4286
$ws.Range("K132").Value = 5569.2858
$ws.Range("M132").Value = -3039.2858

$ws.Range("H136").Value = 3120.8
$ws.Range("I136").Value = 2912.111
$ws.Range("K136").Value = 8736.332999999999
$ws.Range("M136").Value = -6186.332999999999
